$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell B1 from "Número" to "Numero"
$ws.Range("B1").Value = "Numero"

# Update the active cell selection to B2
$ws.Range("B2").Select()
